$wb = $excel.ActiveWorkbook

# Sheets, by name, as referenced in the workbook
$wsProduction = $wb.Worksheets.Item("productionAssets")
$wsStorage    = $wb.Worksheets.Item("storageAssets")

# --- storageAssets: new "District_Battery_500_kWh" storage asset row (row 15) ---
$wsStorage.Range("A15").Value = 14
$wsStorage.Range("B15").Value = "District_Battery_500_kWh"
$wsStorage.Range("C15").Value = "STORAGE"
$wsStorage.Range("D15").Value = "STORAGE_ELECTRIC"
$wsStorage.Range("E15").Value = 100
$wsStorage.Range("F15").Value = 0
$wsStorage.Range("G15").Value = 1
$wsStorage.Range("H15").Value = 0
$wsStorage.Range("I15").Value = 0
$wsStorage.Range("J15").Value = 0
$wsStorage.Range("K15").Value = 0
$wsStorage.Range("L15").Value = 500
$wsStorage.Range("M15").Value = 0
$wsStorage.Range("M15").NumberFormat = "0.00E+00"

# Existing "House_battery" row: stateOfCharge_r 0 -> 1
$wsStorage.Range("G2").Value = 1

# --- Sheet selection / active-view bookkeeping ---
# productionAssets (previously active) keeps its own selection, now at C11
$wsProduction.Activate() | Out-Null
$wsProduction.Range("C11").Select() | Out-Null

# storageAssets becomes the active tab/sheet, selection moves to B19
$wsStorage.Activate() | Out-Null
$wsStorage.Range("B19").Select() | Out-Null
